$wb = $excel.ActiveWorkbook

# --- 1) Rename header labels to match the new PO-forecast naming scheme ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after "Monthly Trend" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, bordered, centered) used on the other sheets
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows
$data = @(
  @(45137.99999999999,321,-26.37485310660766,691.3503931634493),
  @(45144.99999999999,319,-41.13173758998624,697.4988738301176),
  @(45165.99999999999,313,-27.62994863236038,659.1940251715606),
  @(45172.99999999999,311,-56.0001510329251,683.9992030580478),
  @(45179.99999999999,308,-66.42055448565183,661.2544959198401),
  @(45186.99999999999,306,-47.92403700987794,678.2425153290857),
  @(45193.99999999999,304,-54.76496944030786,644.3642098975768),
  @(45200.99999999999,302,-49.46691241018775,683.1452475654792),
  @(45207.99999999999,300,-39.16518460653634,648.9871104096633),
  @(45228.99999999999,294,-61.7462236959167,645.3150905303454),
  @(45235.99999999999,292,-75.47788810933093,632.0845954422282),
  @(45242.99999999999,290,-50.98087888754146,657.9531823684994),
  @(45249.99999999999,288,-102.0754387676903,637.1436314636308),
  @(45270.99999999999,281,-93.53243120583049,654.4510661944362),
  @(45277.99999999999,279,-65.47565449771258,603.8672368109977),
  @(45312.99999999999,269,-62.47940622114317,629.955922032196),
  @(45326.99999999999,265,-86.47130680370653,618.5229746543392),
  @(45333.99999999999,263,-86.34028785995719,605.5015976330302),
  @(45340.99999999999,260,-80.15217669554488,632.7848579332637),
  @(45347.99999999999,258,-72.12166152427066,611.6630978317992),
  @(45354.99999999999,256,-95.79475693486428,650.0727092897188),
  @(45361.99999999999,254,-101.0513602630947,616.0632903432503),
  @(45368.99999999999,252,-103.1495404701203,590.0663789903448),
  @(45382.99999999999,248,-77.98733021764353,590.5921084978763),
  @(45494.99999999999,215,-111.7200621479483,557.5104319915403),
  @(45501.99999999999,213,-119.0014788446964,564.280563854149),
  @(45508.99999999999,210,-150.1178716423307,571.7798991477364),
  @(45515.99999999999,208,-150.0045151050702,573.3691184464665),
  @(45522.99999999999,206,-150.0448210647493,569.6154398191526),
  @(45529.99999999999,204,-152.1438885305891,528.403840372429),
  @(45613.99999999999,179,-174.8209820732615,516.332695241665),
  @(45620.99999999999,177,-186.2753292317971,509.3337329319882),
  @(45627.99999999999,175,-180.6586291431469,539.2268985082916),
  @(45634.99999999999,173,-190.46361901817,530.7310526627302),
  @(45641.99999999999,171,-189.9734868340217,533.1683060389925),
  @(45648.99999999999,169,-173.110341010167,540.1715065313767),
  @(45655.99999999999,167,-195.1974367309816,525.963150058896),
  @(45662.99999999999,165,-166.1893815189382,520.0059449302208),
  @(45669.99999999999,162,-198.7182104334266,507.6495340618631)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $wsForecast.Cells.Item($row, 1).Value = $data[$i][0]
  $wsForecast.Cells.Item($row, 2).Value = $data[$i][1]
  $wsForecast.Cells.Item($row, 3).Value = $data[$i][2]
  $wsForecast.Cells.Item($row, 4).Value = $data[$i][3]
}

# Match the date-formatted style used for column A on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A40").PasteSpecial(-4122)
